$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.54819999999999
$ws.Range("D3").Value = -6.704599999999992
$ws.Range("E8").Value = 16.0849
$ws.Range("E11").Value = 16.5722
$ws.Range("A12").Value = -21.59099999999999
$ws.Range("C14").Value = -13.9523
$ws.Range("E14").Value = 16.8149
$ws.Range("E15").Value = 16.4125
$ws.Range("E17").Value = 16.69600000000001
$ws.Range("D20").Value = -7.585600000000004
$ws.Range("D25").Value = -7.597600000000003
$ws.Range("C26").Value = -13.1398
$ws.Range("E26").Value = 15.80359999999999
$ws.Range("A27").Value = -21.84789999999999
$ws.Range("D30").Value = -7.275200000000007
$ws.Range("C31").Value = -12.8312
$ws.Range("A32").Value = -21.4927
$ws.Range("C35").Value = -12.50160000000001
$ws.Range("A36").Value = -19.82039999999999
$ws.Range("E36").Value = 16.24540000000001
$ws.Range("C37").Value = -14.0817
$ws.Range("A38").Value = -19.4375
$ws.Range("D44").Value = -7.192100000000006
$ws.Range("C45").Value = -14.28149999999999
$ws.Range("A46").Value = -21.4187
$ws.Range("D47").Value = -7.382000000000001
$ws.Range("C52").Value = -10.9667
$ws.Range("A54").Value = -21.38639999999998
$ws.Range("A55").Value = -22.52580000000001
$ws.Range("A56").Value = -22.14960000000001
$ws.Range("C57").Value = -14.0301
$ws.Range("D58").Value = -7.8933
$ws.Range("E64").Value = 17.38200000000001
$ws.Range("A67").Value = -21.46429999999998
$ws.Range("A69").Value = -21.64439999999998
$ws.Range("A72").Value = -21.42709999999998
$ws.Range("D78").Value = -7.699100000000002
$ws.Range("E79").Value = 17.83040000000002
$ws.Range("C81").Value = -13.2058
$ws.Range("A83").Value = -21.39989999999999
$ws.Range("C83").Value = -11.44260000000001
$ws.Range("D84").Value = -8.610800000000006
$ws.Range("A86").Value = -22.33250000000001
$ws.Range("D89").Value = -7.060299999999994
$ws.Range("E89").Value = 17.46360000000002
$ws.Range("A91").Value = -21.5037
$ws.Range("D91").Value = -6.461499999999996
$ws.Range("D92").Value = -6.526
$ws.Range("A93").Value = -21.1114
$ws.Range("D96").Value = -7.457800000000004
$ws.Range("A99").Value = -20.40319999999999
$ws.Range("C100").Value = -12.4678
$ws.Range("C102").Value = -14.70619999999999
$ws.Range("D102").Value = -7.796499999999998
